$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsCodes = $wb.Worksheets.Item("Include from Split Method Cod")

# URL
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/split-method"

# Version
$wsMeta.Range("B3").Value = "8.0.0"

# Date
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# System URI
$wsCodes.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/split-method"
